$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the sheet; this pushes all existing
# data (rows 1-86) down by one row (to rows 2-87) and keeps their values
# intact.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Cells.Item(1, 1).Value = "id"
$ws.Cells.Item(1, 2).Value = "AgendaHTMLiFrameURL"
$ws.Cells.Item(1, 3).Value = "PageURL"
